$wb = $excel.ActiveWorkbook

# Overview sheet: row 3 corresponds to 75b38656-db3a-4516-9551-dfd6146543e2
$wsOverview = $wb.Sheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-12 08:58:35"

# zh-cn sheet: row 3 corresponds to 75b38656-db3a-4516-9551-dfd6146543e2
$wsZhCn = $wb.Sheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-08-12 08:58:28"
$wsZhCn.Range("K3").Value = "2016-08-12 08:58:55"

# de-de sheet: row 3 corresponds to 75b38656-db3a-4516-9551-dfd6146543e2
$wsDeDe = $wb.Sheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-08-12 08:58:35"
$wsDeDe.Range("K3").Value = "2016-08-12 08:59:09"
